# Apply "small fix for Dry Clap 1" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "presets"
$ws.Name = "presets"

# Update cell B40: change from "Dry Clap 1" to a new duplicate string "Dry Clap 1 (dup)"
$ws.Range("B40").Value = "Dry Clap 1 (dup)"

# Update the selection shown in the sheet view to B40
$ws.Range("B40").Select()
